# Naresh Mail id added in Config 30/08/2022 3.27 pm
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting")

# To_MailSend (B3): append naresh.kumar@e5.ai to the end of the cc list
$ws.Range("B3").Value = " mvprasanth97@gmail.com ; lakshmi.u@tiliconveli.com ; narenbagavathye5@gmail.com ; sornalakshmie5@gmail.com ; aartiak.e5@gmail.com ; sankaravenie5@gmail.com ;sharongiftae5@gmail.com ;naresh.kumar@e5.ai"

# B13: drop alagappan.m@e5.ai, keep just naresh.kumar@e5.ai
$ws.Range("B13").Value = "naresh.kumar@e5.ai"

# B14: drop ";alagappan.m@e5.ai" from the middle of the list
$ws.Range("B14").Value = "naresh.kumar@e5.ai ; mvprasanth97@gmail.com ; lakshmi.u@tiliconveli.com ; narenbagavathye5@gmail.com ; sornalakshmie5@gmail.com ; aartiak.e5@gmail.com ; sankaravenie5@gmail.com ;sharongiftae5@gmail.com"

# Update the saved selection/active cell on the Setting sheet
$ws.Range("E16").Select()
